$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Add the new "Vision / back_cam / Vision Sensor" row for PORT 2 (row 3)
$ws.Range("B3").Value = "Vision"
$ws.Range("C3").Value = "back_cam"
$ws.Range("D3").Value = "Vision Sensor"

# Add the new "Limit / rear_switch" entry for PORT G (row 28)
$ws.Range("B28").Value = "Limit"
$ws.Range("C28").Value = "rear_switch"

# Update the view: scroll position and active selection
$ws.Range("B29").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 130
